$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.56887376080371
$ws.Range("C2").Value = 9.616962402972069
$ws.Range("E2").Value = 11.47552508142522
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 24.31200988536654
$ws.Range("H2").Value = 12.79737398264698
$ws.Range("I2").Value = 17.67094968774608
$ws.Range("M2").Value = 14.61234749703238
$ws.Range("B3").Value = 11.910159495426
$ws.Range("C3").Value = 9.059986686841064
$ws.Range("E3").Value = 11.39582216632876
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 24.21090572750521
$ws.Range("H3").Value = 12.86216391753769
$ws.Range("I3").Value = 17.82133926884983
$ws.Range("M3").Value = 14.30045080280758
$ws.Range("B4").Value = 11.48745990589954
$ws.Range("C4").Value = 8.698360605248995
$ws.Range("E4").Value = 11.35117858511622
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 24.16648425500686
$ws.Range("H4").Value = 12.90609710775776
$ws.Range("I4").Value = 17.92035240763592
$ws.Range("M4").Value = 14.10899837751401
$ws.Range("B5").Value = 11.31079956168846
$ws.Range("C5").Value = 8.546097223116046
$ws.Range("E5").Value = 11.33407995066704
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 24.15281382906065
$ws.Range("H5").Value = 12.92503774980555
$ws.Range("I5").Value = 17.96236948673096
$ws.Range("M5").Value = 14.03110750795511
$ws.Range("B6").Value = 11.28120504366204
$ws.Range("C6").Value = 8.520519614600662
$ws.Range("E6").Value = 11.33130717530942
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 24.15081105509914
$ws.Range("H6").Value = 12.92824530054962
$ws.Range("I6").Value = 17.96944688279601
$ws.Range("M6").Value = 14.01818487443223
$ws.Range("B7").Value = 11.48509498173897
$ws.Range("C7").Value = 8.69632689654027
$ws.Range("E7").Value = 11.35094354056625
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 24.16628196696071
$ws.Range("H7").Value = 12.90634835575482
$ws.Range("I7").Value = 17.92091232191434
$ws.Range("M7").Value = 14.10794724499572
$ws.Range("B8").Value = 12.34563501984533
$ws.Range("C8").Value = 9.429014812300268
$ws.Range("E8").Value = 11.44716144666769
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 24.27348315885691
$ws.Range("H8").Value = 12.81884791711931
$ws.Range("I8").Value = 17.72141308964441
$ws.Range("M8").Value = 14.50486206528624
$ws.Range("B9").Value = 13.8819435843696
$ws.Range("C9").Value = 10.70869662635822
$ws.Range("E9").Value = 11.66918175842961
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 24.62357552597775
$ws.Range("H9").Value = 12.68049270423171
$ws.Range("I9").Value = 17.38358854304662
$ws.Range("M9").Value = 15.27877483991484
$ws.Range("B10").Value = 14.91167324054975
$ws.Range("C10").Value = 11.55197424508465
$ws.Range("E10").Value = 11.85148439154785
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 24.96490727222377
$ws.Range("H10").Value = 12.59949499457067
$ws.Range("I10").Value = 17.1685639540186
$ws.Range("M10").Value = 15.83843880154381
$ws.Range("B11").Value = 15.35760911795355
$ws.Range("C11").Value = 11.91445189657436
$ws.Range("E11").Value = 11.93830104906722
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 25.13796590495555
$ws.Range("H11").Value = 12.56721402892345
$ws.Range("I11").Value = 17.07808887489958
$ws.Range("M11").Value = 16.08988543424344
$ws.Range("B12").Value = 15.52317536274853
$ws.Range("C12").Value = 12.04866989175218
$ws.Range("E12").Value = 11.97170890473455
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 25.20599973536239
$ws.Range("H12").Value = 12.55565323583752
$ws.Range("I12").Value = 17.04489628981304
$ws.Range("M12").Value = 16.18455512097112
$ws.Range("B13").Value = 15.48766537924742
$ws.Range("C13").Value = 12.01989909502442
$ws.Range("E13").Value = 11.96449068939658
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 25.19123724511072
$ws.Range("H13").Value = 12.55811345120938
$ws.Range("I13").Value = 17.05199718343609
$ws.Range("M13").Value = 16.16419211840904
$ws.Range("B14").Value = 15.37129681700192
$ws.Range("C14").Value = 11.92555512493025
$ws.Range("E14").Value = 11.94103902283489
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 25.1435133417218
$ws.Range("H14").Value = 12.5662495831561
$ws.Range("I14").Value = 17.07533661475297
$ws.Range("M14").Value = 16.09768546607004
$ws.Range("B15").Value = 15.29958624173079
$ws.Range("C15").Value = 11.86737014588224
$ws.Range("E15").Value = 11.92674270606538
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 25.1146048202573
$ws.Range("H15").Value = 12.57131978178306
$ws.Range("I15").Value = 17.08977220883328
$ws.Range("M15").Value = 16.05687417629974
$ws.Range("B16").Value = 14.88207174851346
$ws.Range("C16").Value = 11.5278596660187
$ws.Range("E16").Value = 11.84588668056121
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 24.95395070083527
$ws.Range("H16").Value = 12.60169708287471
$ws.Range("I16").Value = 17.17462551036413
$ws.Range("M16").Value = 15.82193457120888
$ws.Range("B17").Value = 14.62012803111095
$ws.Range("C17").Value = 11.3141630855048
$ws.Range("E17").Value = 11.79726103490838
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 24.85991225641118
$ws.Range("H17").Value = 12.62150674483756
$ws.Range("I17").Value = 17.22856935959316
$ws.Range("M17").Value = 15.67693103450647
$ws.Range("B18").Value = 14.46735071195719
$ws.Range("C18").Value = 11.18926133225462
$ws.Range("E18").Value = 11.76966033009945
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 24.80750193116981
$ws.Range("H18").Value = 12.6333300444563
$ws.Range("I18").Value = 17.26028658650515
$ws.Range("M18").Value = 15.59323790623805
$ws.Range("B19").Value = 14.41526203912447
$ws.Range("C19").Value = 11.14663033715737
$ws.Range("E19").Value = 11.76037910770588
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 24.79004649142204
$ws.Range("H19").Value = 12.63740673452426
$ws.Range("I19").Value = 17.27114364553805
$ws.Range("M19").Value = 15.56485388126776
$ws.Range("B20").Value = 14.64823171827339
$ws.Range("C20").Value = 11.33711736990056
$ws.Range("E20").Value = 11.80239947325299
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 24.86974951577688
$ws.Range("H20").Value = 12.61935349133508
$ws.Range("I20").Value = 17.22275542556486
$ws.Range("M20").Value = 15.69239772695204
$ws.Range("B21").Value = 15.40556710113534
$ws.Range("C21").Value = 11.95334888611119
$ws.Range("E21").Value = 11.94791311839159
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 25.15746365578941
$ws.Range("H21").Value = 12.56384174816487
$ws.Range("I21").Value = 17.06845215912991
$ws.Range("M21").Value = 16.11723567403264
$ws.Range("B22").Value = 15.88127326530437
$ws.Range("C22").Value = 12.33834646166475
$ws.Range("E22").Value = 12.04610397748409
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 25.36004188755766
$ws.Range("H22").Value = 12.53143096491612
$ws.Range("I22").Value = 16.97383950345283
$ws.Range("M22").Value = 16.39165966674587
$ws.Range("B23").Value = 15.62916016417399
$ws.Range("C23").Value = 12.13449060198183
$ws.Range("E23").Value = 11.99342410673993
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 25.25061297146399
$ws.Range("H23").Value = 12.54837293960811
$ws.Range("I23").Value = 17.02376137356202
$ws.Range("M23").Value = 16.24551999702043
$ws.Range("B24").Value = 14.63553283386399
$ws.Range("C24").Value = 11.32674610647375
$ws.Range("E24").Value = 11.80007527805839
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 24.86529693576183
$ws.Range("H24").Value = 12.62032562443806
$ws.Range("I24").Value = 17.2253817125925
$ws.Range("M24").Value = 15.68540625172178
$ws.Range("B25").Value = 13.4833102169623
$ws.Range("C25").Value = 10.37947153275012
$ws.Range("E25").Value = 11.60565993429262
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 24.51396487365728
$ws.Range("H25").Value = 12.71432288758886
$ws.Range("I25").Value = 17.46919973658539
$ws.Range("M25").Value = 15.07054676107623
